$d = $word.ActiveDocument
$t = $d.Tables(1)

# Update existing cell text for rows 1-16 and rows 17-19 (pre-insertion numbering)
$t.Cell(1,1).Range.Text = "72-26="
$t.Cell(1,2).Range.Text = "60-31="
$t.Cell(1,3).Range.Text = "63-27="
$t.Cell(1,4).Range.Text = "32+29="
$t.Cell(1,5).Range.Text = "41-16="
$t.Cell(2,1).Range.Text = "93-85="
$t.Cell(2,2).Range.Text = "66-47="
$t.Cell(2,3).Range.Text = "58-19="
$t.Cell(2,4).Range.Text = "53+29="
$t.Cell(2,5).Range.Text = "70-56="
$t.Cell(3,1).Range.Text = "95-39="
$t.Cell(3,2).Range.Text = "56-18="
$t.Cell(3,3).Range.Text = "92-48="
$t.Cell(3,4).Range.Text = "95-69="
$t.Cell(3,5).Range.Text = "36-8="
$t.Cell(4,1).Range.Text = "9+72="
$t.Cell(4,2).Range.Text = "76-67="
$t.Cell(4,3).Range.Text = "56+39="
$t.Cell(4,4).Range.Text = "35+36="
$t.Cell(4,5).Range.Text = "28+27="
$t.Cell(5,1).Range.Text = "41-6="
$t.Cell(5,2).Range.Text = "80-11="
$t.Cell(5,3).Range.Text = "44+47="
$t.Cell(5,4).Range.Text = "95-27="
$t.Cell(5,5).Range.Text = "34+57="
$t.Cell(6,1).Range.Text = "17+15="
$t.Cell(6,2).Range.Text = "50-2="
$t.Cell(6,3).Range.Text = "24-7="
$t.Cell(6,4).Range.Text = "45-9="
$t.Cell(6,5).Range.Text = "91-15="
$t.Cell(7,1).Range.Text = "53+9="
$t.Cell(7,2).Range.Text = "65+18="
$t.Cell(7,3).Range.Text = "22-6="
$t.Cell(7,4).Range.Text = "90-65="
$t.Cell(7,5).Range.Text = "40-14="
$t.Cell(8,1).Range.Text = "16+26="
$t.Cell(8,2).Range.Text = "29+63="
$t.Cell(8,3).Range.Text = "57+27="
$t.Cell(8,4).Range.Text = "63-59="
$t.Cell(8,5).Range.Text = "48+18="
$t.Cell(9,1).Range.Text = "57+39="
$t.Cell(9,2).Range.Text = "71-49="
$t.Cell(9,3).Range.Text = "90-69="
$t.Cell(9,4).Range.Text = "65-58="
$t.Cell(9,5).Range.Text = "80-53="
$t.Cell(10,1).Range.Text = "19+49="
$t.Cell(10,2).Range.Text = "94-85="
$t.Cell(10,3).Range.Text = "83-37="
$t.Cell(10,4).Range.Text = "14+38="
$t.Cell(10,5).Range.Text = "32+49="
$t.Cell(11,1).Range.Text = "56+38="
$t.Cell(11,2).Range.Text = "34-8="
$t.Cell(11,3).Range.Text = "59+14="
$t.Cell(11,4).Range.Text = "5+89="
$t.Cell(11,5).Range.Text = "69+27="
$t.Cell(12,1).Range.Text = "23+19="
$t.Cell(12,2).Range.Text = "70-18="
$t.Cell(12,3).Range.Text = "84-5="
$t.Cell(12,4).Range.Text = "61-13="
$t.Cell(12,5).Range.Text = "39+49="
$t.Cell(13,1).Range.Text = "14+79="
$t.Cell(13,2).Range.Text = "47+48="
$t.Cell(13,3).Range.Text = "16+19="
$t.Cell(13,4).Range.Text = "56+5="
$t.Cell(13,5).Range.Text = "20-12="
$t.Cell(14,1).Range.Text = "19+23="
$t.Cell(14,2).Range.Text = "58-9="
$t.Cell(14,3).Range.Text = "94-56="
$t.Cell(14,4).Range.Text = "80-23="
$t.Cell(14,5).Range.Text = "34-9="
$t.Cell(15,1).Range.Text = "81-18="
$t.Cell(15,2).Range.Text = "91-73="
$t.Cell(15,3).Range.Text = "12+49="
$t.Cell(15,4).Range.Text = "90-38="
$t.Cell(15,5).Range.Text = "70-11="
$t.Cell(16,1).Range.Text = "97-89="
$t.Cell(16,2).Range.Text = "66+19="
$t.Cell(16,3).Range.Text = "90-28="
$t.Cell(16,4).Range.Text = "65-36="
$t.Cell(16,5).Range.Text = "91-84="
$t.Cell(17,2).Range.Text = "94-27="
$t.Cell(17,3).Range.Text = "80-6="
$t.Cell(17,4).Range.Text = "77+14="
$t.Cell(17,5).Range.Text = "15+59="
$t.Cell(18,1).Range.Text = "94-57="
$t.Cell(18,2).Range.Text = "23+48="
$t.Cell(18,3).Range.Text = "40-19="
$t.Cell(18,4).Range.Text = "78-9="
$t.Cell(18,5).Range.Text = "26+39="
$t.Cell(19,1).Range.Text = "37+17="
$t.Cell(19,2).Range.Text = "80-69="
$t.Cell(19,3).Range.Text = "70-18="
$t.Cell(19,4).Range.Text = "72-63="
$t.Cell(19,5).Range.Text = "77+5="

# Insert new row before (original) row 17, which becomes the new row 17
$beforeRow = $t.Rows(17)
$newRow = $t.Rows.Add($beforeRow)
$newRowValues = @("29+33=", "75+8=", "37+14=", "34-18=", "60-16=")
for ($c = 1; $c -le 5; $c++) {
    $t.Cell(17, $c).Range.Text = $newRowValues[$c-1]
}

# Delete the last row (originally row 20, now row 21 after insertion)
$t.Rows.Last.Delete()

Write-Output "done"